{"js": "// Update each \"A\u00d7B=C\" answer cell in the practice-sheet table to the\n// newly generated equation/answer. Every old string below occurs exactly\n// once in the document, so a plain body.search + insertText replace is\n// sufficient (no need to address cells by row/column index).\nconst replacements = [\n  [\"31\u00d776=2356\", \"79\u00d737=2923\"],\n  [\"69\u00d766=4554\", \"59\u00d774=4366\"],\n  [\"57\u00d717=969\", \"48\u00d734=1632\"],\n  [\"50\u00d725=1250\", \"74\u00d714=1036\"],\n  [\"60\u00d719=1140\", \"20\u00d782=1640\"],\n  [\"47\u00d745=2115\", \"47\u00d739=1833\"],\n  [\"15\u00d717=255\", \"85\u00d757=4845\"],\n  [\"32\u00d772=2304\", \"63\u00d736=2268\"],\n  [\"31\u00d794=2914\", \"87\u00d763=5481\"],\n  [\"17\u00d764=1088\", \"94\u00d796=9024\"],\n  [\"35\u00d799=3465\", \"50\u00d751=2550\"],\n  [\"16\u00d770=1120\", \"28\u00d723=644\"],\n  [\"25\u00d773=1825\", \"70\u00d768=4760\"],\n  [\"83\u00d781=6723\", \"32\u00d713=416\"],\n  [\"39\u00d789=3471\", \"37\u00d752=1924\"],\n  [\"54\u00d786=4644\", \"42\u00d740=1680\"],\n  [\"51\u00d795=4845\", \"85\u00d730=2550\"],\n  [\"39\u00d768=2652\", \"70\u00d738=2660\"],\n  [\"89\u00d732=2848\", \"44\u00d773=3212\"],\n  [\"97\u00d713=1261\", \"27\u00d718=486\"],\n  [\"99\u00d737=3663\", \"94\u00d734=3196\"],\n  [\"55\u00d794=5170\", \"80\u00d771=5680\"],\n  [\"88\u00d774=6512\", \"29\u00d748=1392\"],\n  [\"33\u00d760=1980\", \"68\u00d773=4964\"],\n  [\"82\u00d778=6396\", \"86\u00d763=5418\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n", "ps1": "# Update each \"A\u00d7B=C\" answer cell in the practice-sheet table to the\n# newly generated equation/answer. Every old string below occurs exactly\n# once in the document, so Find/Replace over the whole document body is\n# sufficient (no need to address cells by row/column index).\n$d = $word.ActiveDocument\n$pairs = @(\n  @(\"31\u00d776=2356\", \"79\u00d737=2923\"),\n  @(\"69\u00d766=4554\", \"59\u00d774=4366\"),\n  @(\"57\u00d717=969\", \"48\u00d734=1632\"),\n  @(\"50\u00d725=1250\", \"74\u00d714=1036\"),\n  @(\"60\u00d719=1140\", \"20\u00d782=1640\"),\n  @(\"47\u00d745=2115\", \"47\u00d739=1833\"),\n  @(\"15\u00d717=255\", \"85\u00d757=4845\"),\n  @(\"32\u00d772=2304\", \"63\u00d736=2268\"),\n  @(\"31\u00d794=2914\", \"87\u00d763=5481\"),\n  @(\"17\u00d764=1088\", \"94\u00d796=9024\"),\n  @(\"35\u00d799=3465\", \"50\u00d751=2550\"),\n  @(\"16\u00d770=1120\", \"28\u00d723=644\"),\n  @(\"25\u00d773=1825\", \"70\u00d768=4760\"),\n  @(\"83\u00d781=6723\", \"32\u00d713=416\"),\n  @(\"39\u00d789=3471\", \"37\u00d752=1924\"),\n  @(\"54\u00d786=4644\", \"42\u00d740=1680\"),\n  @(\"51\u00d795=4845\", \"85\u00d730=2550\"),\n  @(\"39\u00d768=2652\", \"70\u00d738=2660\"),\n  @(\"89\u00d732=2848\", \"44\u00d773=3212\"),\n  @(\"97\u00d713=1261\", \"27\u00d718=486\"),\n  @(\"99\u00d737=3663\", \"94\u00d734=3196\"),\n  @(\"55\u00d794=5170\", \"80\u00d771=5680\"),\n  @(\"88\u00d774=6512\", \"29\u00d748=1392\"),\n  @(\"33\u00d760=1980\", \"68\u00d773=4964\"),\n  @(\"82\u00d778=6396\", \"86\u00d763=5418\"),\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $found = $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n    if (-not $found) {\n        throw \"No match found for: $old\"\n    }\n}\n\n"}
